$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $s
}

$ws.Cells.Item(2, 4).Value = '23.131.93'
$ws.Cells.Item(2, 5).Value = '  -3.09%  '

$ws.Cells.Item(3, 4).Value = '1.606.66'
$ws.Cells.Item(3, 5).Value = '  -2.45%  '

$ws.Cells.Item(4, 5).Value = '  -0.11%  '

Set-TextValue 5 4 '1.002'
$ws.Cells.Item(5, 5).Value = '  +0.00%  '

Set-TextValue 6 4 '301.94'
$ws.Cells.Item(6, 5).Value = '  -2.22%  '

Set-TextValue 7 4 '0.3765'
$ws.Cells.Item(7, 5).Value = '  -3.15%  '

Set-TextValue 8 4 '0.3632'
$ws.Cells.Item(8, 5).Value = '  -5.06%  '

Set-TextValue 9 4 '48.62'
$ws.Cells.Item(9, 5).Value = '  -5.40%  '

$ws.Cells.Item(10, 5).Value = '  -0.12%  '

Set-TextValue 11 4 '1.260'
$ws.Cells.Item(11, 5).Value = '  -6.37%  '

Set-TextValue 12 4 '0.08045'
$ws.Cells.Item(12, 5).Value = '  -4.50%  '

Set-TextValue 13 4 '22.77'
$ws.Cells.Item(13, 5).Value = '  -4.50%  '

Set-TextValue 14 4 '6.549'
$ws.Cells.Item(14, 5).Value = '  -7.33%  '

Set-TextValue 15 4 '7.471'
$ws.Cells.Item(15, 5).Value = '  -5.31%  '

Set-TextValue 16 4 '0.00001253'
$ws.Cells.Item(16, 5).Value = '  -4.70%  '

$ws.Cells.Item(17, 4).Value = '1.603.68'
$ws.Cells.Item(17, 5).Value = '  -2.88%  '

Set-TextValue 18 4 '91.28'
$ws.Cells.Item(18, 5).Value = '  -3.19%  '

$ws.Cells.Item(19, 5).Value = '  -2.82%  '

Set-TextValue 20 4 '18.25'
$ws.Cells.Item(20, 5).Value = '  -7.02%  '

Set-TextValue 21 4 '6.542'
$ws.Cells.Item(21, 5).Value = '  -5.47%  '

Set-TextValue 22 4 '1.002'
$ws.Cells.Item(22, 5).Value = '  +0.01%  '

Set-TextValue 23 4 '13.02'
$ws.Cells.Item(23, 5).Value = '  -4.65%  '

$ws.Cells.Item(24, 4).Value = '23.157.42'
$ws.Cells.Item(24, 5).Value = '  -3.02%  '

$ws.Cells.Item(25, 5).Value = '  -3.73%  '

Set-TextValue 26 4 '2.835'
$ws.Cells.Item(26, 5).Value = '  -4.40%  '

Set-TextValue 27 4 '21.02'
$ws.Cells.Item(27, 5).Value = '  -4.40%  '

Set-TextValue 28 4 '150.26'
$ws.Cells.Item(28, 5).Value = '  -0.30%  '

Set-TextValue 29 4 '5.261'
$ws.Cells.Item(29, 5).Value = '  -2.25%  '

Set-TextValue 30 4 '131.44'
$ws.Cells.Item(30, 5).Value = '  -5.01%  '

Set-TextValue 31 4 '2.397'
$ws.Cells.Item(31, 5).Value = '  -4.45%  '

Set-TextValue 32 4 '6.744'
$ws.Cells.Item(32, 5).Value = '  -13.26%  '

$ws.Cells.Item(33, 4).Value = '1.778.67'
$ws.Cells.Item(33, 5).Value = '  -2.89%  '

Set-TextValue 34 4 '0.9629'
$ws.Cells.Item(34, 5).Value = '  -8.43%  '

Set-TextValue 35 4 '0.07693'

Set-TextValue 36 4 '0.02766'
$ws.Cells.Item(36, 5).Value = '  -6.35%  '

Set-TextValue 37 4 '0.2541'
$ws.Cells.Item(37, 5).Value = '  -5.03%  '

Set-TextValue 38 4 '6.190'
$ws.Cells.Item(38, 5).Value = '  -7.84%  '

Set-TextValue 39 4 '10.10'
$ws.Cells.Item(39, 5).Value = '  -6.87%  '

Set-TextValue 40 4 '0.08831'
$ws.Cells.Item(40, 5).Value = '  -2.91%  '

Set-TextValue 41 4 '1.390'
$ws.Cells.Item(41, 5).Value = '  -2.16%  '

Set-TextValue 42 4 '0.7146'
$ws.Cells.Item(42, 5).Value = '  -5.45%  '

Set-TextValue 43 4 '12.69'
$ws.Cells.Item(43, 5).Value = '  -5.59%  '

Set-TextValue 44 4 '15.81'
$ws.Cells.Item(44, 5).Value = '  -3.62%  '

$ws.Cells.Item(45, 5).Value = '  -4.87%  '

Set-TextValue 46 4 '1.001'
$ws.Cells.Item(46, 5).Value = '  +0.16%  '

Set-TextValue 47 4 '2.287'
$ws.Cells.Item(47, 5).Value = '  -6.77%  '

Set-TextValue 48 4 '3.973'
$ws.Cells.Item(48, 5).Value = '  -2.64%  '

$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 49 4 '131.79'
$ws.Cells.Item(49, 5).Value = '  -1.57%  '

$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 50 4 '0.07976'
$ws.Cells.Item(50, 5).Value = '  -3.58%  '

Set-TextValue 51 4 '1.168'
$ws.Cells.Item(51, 5).Value = '  -3.17%  '
